$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New column F: "president" header (same formatting as the other headers)
# plus "Ford" in every data row.
$ws.Range("A1").Copy()
$ws.Range("F1").PasteSpecial(-4122)   # xlPasteFormats - copy header style (bold/centered/bordered) onto F1
$ws.Range("F1").Value = "president"

$ws.Range("F2:F37").Value = "Ford"
